$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050528638925254
$ws.Range("D2").Value = 1.048496406949045
$ws.Range("E2").Value = 1.057278510796025
$ws.Range("F2").Value = 1.067571630651033
$ws.Range("I2").Value = 1.036456859639565
$ws.Range("J2").Value = 1.055561314756512
$ws.Range("K2").Value = 1.051256089456256
$ws.Range("L2").Value = 1.06001391771789
$ws.Range("M2").Value = 1.070279135350954
$ws.Range("N2").Value = 1.057060332898042
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052277394149757
$ws.Range("D3").Value = 1.049812534128018
$ws.Range("E3").Value = 1.058911033885642
$ws.Range("F3").Value = 1.069462546953422
$ws.Range("I3").Value = 1.036827886014476
$ws.Range("J3").Value = 1.056955887792743
$ws.Range("K3").Value = 1.052382969914302
$ws.Range("L3").Value = 1.061458166361066
$ws.Range("M3").Value = 1.071983189389152
$ws.Range("N3").Value = 1.058456886387943
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053405986869099
$ws.Range("D4").Value = 1.050661316087223
$ws.Range("E4").Value = 1.059964849921135
$ws.Range("F4").Value = 1.07068382793268
$ws.Range("I4").Value = 1.037065383144963
$ws.Range("J4").Value = 1.057855016510148
$ws.Range("K4").Value = 1.053108760611559
$ws.Range("L4").Value = 1.062389683337169
$ws.Range("M4").Value = 1.073083133202091
$ws.Range("N4").Value = 1.059357291971259
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053879751762024
$ws.Range("D5").Value = 1.05101747450394
$ws.Range("E5").Value = 1.060407280636759
$ws.Range("F5").Value = 1.071196727175035
$ws.Range("I5").Value = 1.037164612656027
$ws.Range("J5").Value = 1.058232243547891
$ws.Range("K5").Value = 1.053413084552657
$ws.Range("L5").Value = 1.062780585335306
$ws.Range("M5").Value = 1.073544919360419
$ws.Range("N5").Value = 1.059735054714664
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053959258629429
$ws.Range("D6").Value = 1.051077236110433
$ws.Range("E6").Value = 1.060481532284116
$ws.Range("F6").Value = 1.071282814749908
$ws.Range("I6").Value = 1.037181237815144
$ws.Range("J6").Value = 1.058295537020639
$ws.Range("K6").Value = 1.053464135342626
$ws.Range("L6").Value = 1.062846178359608
$ws.Range("M6").Value = 1.073622418751153
$ws.Range("N6").Value = 1.059798438071404
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053412320050688
$ws.Range("D7").Value = 1.05066607770977
$ws.Range("E7").Value = 1.05997076401577
$ws.Range("F7").Value = 1.070690683367754
$ws.Range("I7").Value = 1.037066711463009
$ws.Range("J7").Value = 1.057860060032106
$ws.Range("K7").Value = 1.053112830126874
$ws.Range("L7").Value = 1.062394909352599
$ws.Range("M7").Value = 1.073089306067807
$ws.Range("N7").Value = 1.059362342655596
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.051120262886017
$ws.Range("D8").Value = 1.048941792728044
$ws.Range("E8").Value = 1.05783076276899
$ws.Range("F8").Value = 1.068211153861196
$ws.Range("I8").Value = 1.036582786601746
$ws.Range("J8").Value = 1.056033298216218
$ws.Range("K8").Value = 1.051637629014216
$ws.Range("L8").Value = 1.060502638786597
$ws.Range("M8").Value = 1.070855594787259
$ws.Range("N8").Value = 1.057532986628398
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047057939490574
$ws.Range("D9").Value = 1.045881156491988
$ws.Range("E9").Value = 1.054039780219328
$ws.Range("F9").Value = 1.06382382339673
$ws.Range("I9").Value = 1.035710104655525
$ws.Range("J9").Value = 1.052788845272026
$ws.Range("K9").Value = 1.049011834042863
$ws.Range("L9").Value = 1.057144611509583
$ws.Range("M9").Value = 1.066898230036181
$ws.Range("N9").Value = 1.054283926188881
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044333007124472
$ws.Range("D10").Value = 1.043825124168709
$ws.Range("E10").Value = 1.051498139732728
$ws.Range("F10").Value = 1.060885781310262
$ws.Range("I10").Value = 1.035114665192634
$ws.Range("J10").Value = 1.050607983374657
$ws.Range("K10").Value = 1.047242995511265
$ws.Range("L10").Value = 1.05488926404603
$ws.Range("M10").Value = 1.064244760289593
$ws.Range("N10").Value = 1.052099967217469
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043148899965227
$ws.Range("D11").Value = 1.042930988199014
$ws.Range("E11").Value = 1.05039399389662
$ws.Range("F11").Value = 1.059610233272208
$ws.Range("I11").Value = 1.034853540878897
$ws.Range("J11").Value = 1.049659229042471
$ws.Range("K11").Value = 1.04647258618044
$ws.Range("L11").Value = 1.05390854409118
$ws.Range("M11").Value = 1.063091961320046
$ws.Range("N11").Value = 1.051149865545308
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042708421075772
$ws.Range("D12").Value = 1.042598273329866
$ws.Range("E12").Value = 1.049983307724162
$ws.Range("F12").Value = 1.059135914350814
$ws.Range("I12").Value = 1.034756047857929
$ws.Range("J12").Value = 1.049306139352295
$ws.Range("K12").Value = 1.046185734998718
$ws.Range("L12").Value = 1.053543623706439
$ws.Range("M12").Value = 1.062663167355153
$ws.Range("N12").Value = 1.050796274427272
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042802934954245
$ws.Range("D13").Value = 1.042669668841372
$ws.Range("E13").Value = 1.050071426787076
$ws.Range("F13").Value = 1.059237681375523
$ws.Range("I13").Value = 1.034776983137258
$ws.Range("J13").Value = 1.049381909307889
$ws.Range("K13").Value = 1.046247296829901
$ws.Range("L13").Value = 1.053621929457208
$ws.Range("M13").Value = 1.062755172243591
$ws.Range("N13").Value = 1.050872151984894
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043112503186733
$ws.Range("D14").Value = 1.042903498055955
$ws.Range("E14").Value = 1.050360057913758
$ws.Range("F14").Value = 1.059571036715777
$ws.Range("I14").Value = 1.034845492299469
$ws.Range("J14").Value = 1.049630056499659
$ws.Range("K14").Value = 1.046448889054975
$ws.Range("L14").Value = 1.053878392749334
$ws.Range("M14").Value = 1.06305652928006
$ws.Range("N14").Value = 1.05112065157414
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043303151859102
$ws.Range("D15").Value = 1.04304748898869
$ws.Range("E15").Value = 1.050537818737031
$ws.Range("F15").Value = 1.0597763579809
$ws.Range("I15").Value = 1.034887636679974
$ws.Range("J15").Value = 1.049782857603712
$ws.Range("K15").Value = 1.046573005296009
$ws.Range("L15").Value = 1.054036323332253
$ws.Range("M15").Value = 1.063242126167995
$ws.Range("N15").Value = 1.051273669673286
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044411502470198
$ws.Range("D16").Value = 1.04388438258997
$ws.Range("E16").Value = 1.051571341023704
$ws.Range("F16").Value = 1.06097036287251
$ws.Range("I16").Value = 1.035131925323415
$ws.Range("J16").Value = 1.050670854530742
$ws.Range("K16").Value = 1.047294029392464
$ws.Range("L16").Value = 1.054954262672493
$ws.Range("M16").Value = 1.064321185569825
$ws.Range("N16").Value = 1.052162927657807
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04510560565049
$ws.Range("D17").Value = 1.044408301184233
$ws.Range("E17").Value = 1.052218667123161
$ws.Range("F17").Value = 1.06171841841493
$ws.Range("I17").Value = 1.035284275593797
$ws.Range("J17").Value = 1.05122667583086
$ws.Range("K17").Value = 1.047745097941659
$ws.Range("L17").Value = 1.055528942662973
$ws.Range("M17").Value = 1.06499701335642
$ws.Range("N17").Value = 1.052719538287925
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045510060961254
$ws.Range("D18").Value = 1.044713522447923
$ws.Range("E18").Value = 1.052595895389186
$ws.Range("F18").Value = 1.062154424026513
$ws.Range("I18").Value = 1.035372821325605
$ws.Range("J18").Value = 1.051550450899088
$ws.Range("K18").Value = 1.048007766151233
$ws.Range("L18").Value = 1.055863745414477
$ws.Range("M18").Value = 1.065390844036131
$ws.Range("N18").Value = 1.053043773153888
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045647901922151
$ws.Range("D19").Value = 1.044817532373476
$ws.Range("E19").Value = 1.052724462274183
$ws.Range("F19").Value = 1.062303036547165
$ws.Range("I19").Value = 1.035402959404299
$ws.Range("J19").Value = 1.051660778009289
$ws.Range("K19").Value = 1.048097256239431
$ws.Range("L19").Value = 1.055977837394225
$ws.Range("M19").Value = 1.065525068197733
$ws.Range("N19").Value = 1.05315425694124
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045031176801869
$ws.Range("D20").Value = 1.044352128206167
$ws.Range("E20").Value = 1.052149251076885
$ws.Range("F20").Value = 1.061638192634996
$ws.Range("I20").Value = 1.035267962737269
$ws.Range("J20").Value = 1.051167085646725
$ws.Range("K20").Value = 1.047696747372895
$ws.Range("L20").Value = 1.055467326254488
$ws.Range("M20").Value = 1.064924541647712
$ws.Range("N20").Value = 1.052659863478893
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043021361118877
$ws.Range("D21").Value = 1.042834657653861
$ws.Range("E21").Value = 1.050275078757734
$ws.Range("F21").Value = 1.059472886486557
$ws.Range("I21").Value = 1.034825331894977
$ws.Range("J21").Value = 1.049557002231862
$ws.Range("K21").Value = 1.046389544228455
$ws.Range("L21").Value = 1.05380288841987
$ws.Range("M21").Value = 1.062967803660582
$ws.Range("N21").Value = 1.051047493560902
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041753947211153
$ws.Range("D22").Value = 1.041877126470264
$ws.Range("E22").Value = 1.049093479823396
$ws.Range("F22").Value = 1.058108433528277
$ws.Range("I22").Value = 1.034544138078294
$ws.Range("J22").Value = 1.048540736535283
$ws.Range("K22").Value = 1.045563673654795
$ws.Range("L22").Value = 1.052752695156321
$ws.Range("M22").Value = 1.061734084354961
$ws.Range("N22").Value = 1.050029784650466
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042426190310532
$ws.Range("D23").Value = 1.042385062049829
$ws.Range("E23").Value = 1.049720179755409
$ws.Range("F23").Value = 1.058832050432328
$ws.Range("I23").Value = 1.034693480190616
$ws.Range("J23").Value = 1.049079856959734
$ws.Range("K23").Value = 1.04600186464725
$ws.Range("L23").Value = 1.053309777918884
$ws.Range("M23").Value = 1.062388434423745
$ws.Range("N23").Value = 1.050569670687759
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045064809240781
$ws.Range("D24").Value = 1.044377511505093
$ws.Range("E24").Value = 1.052180618269417
$ws.Range("F24").Value = 1.06167444420798
$ws.Range("I24").Value = 1.035275334796292
$ws.Range("J24").Value = 1.051194013202303
$ws.Range("K24").Value = 1.047718596249783
$ws.Range("L24").Value = 1.055495169287942
$ws.Range("M24").Value = 1.064957289630926
$ws.Range("N24").Value = 1.052686829274689
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048111022400185
$ws.Range("D25").Value = 1.046675106153855
$ws.Range("E25").Value = 1.055022300693879
$ws.Range("F25").Value = 1.064960298392646
$ws.Range("I25").Value = 1.035938102102679
$ws.Range("J25").Value = 1.053630710816539
$ws.Range("K25").Value = 1.049693845759326
$ws.Range("L25").Value = 1.058015620343417
$ws.Range("M25").Value = 1.067923919482813
$ws.Range("N25").Value = 1.05512698727903

Write-Host "done"